$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Copy the existing centered/wrap-text formatting (style index 1,
#    currently living on B6:B12 / F6:F12) onto the new target ranges
#    before we touch/clear the source cells. Using Copy + PasteSpecial
#    (formats only) re-uses the existing style record instead of
#    minting a new cellXfs entry.
# ------------------------------------------------------------------
$ws.Range("B6").Copy()
$ws.Range("C6:C18").PasteSpecial(-4122)
$ws.Range("B23:B27").PasteSpecial(-4122)
$ws.Range("D23:D27").PasteSpecial(-4122)
$ws.Range("F23:F27").PasteSpecial(-4122)
$ws.Range("B28:B32").PasteSpecial(-4122)
$ws.Range("E28:E32").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2) Write the new course values onto the top-left cell of each future
#    merged block, in the exact order new shared strings are expected:
#    GS-GS-6600 (1:15-2:15), GS-GS-6400 (2:30-3:30), GS-NE-6112 (9-12).
# ------------------------------------------------------------------
$ws.Range("B23").Value = "GS-GS-6600 `n1:15-2:15 `n"
$ws.Range("D23").Value = "GS-GS-6600 `n1:15-2:15 `n"
$ws.Range("F23").Value = "GS-GS-6600 `n1:15-2:15 `n"

$ws.Range("B28").Value = "GS-GS-6400 `n2:30-3:30 `n"
$ws.Range("E28").Value = "GS-GS-6400 `n2:30-3:30 `n"

$ws.Range("C6").Value = "GS-NE-6112 `n09:00-12:00 `nN.0150.01 NRI "

# ------------------------------------------------------------------
# 3) Merge the new course blocks. Order matches the target mergeCells
#    sequence (B23/D23/F23/B28/E28, then C6:C18 last).
# ------------------------------------------------------------------
$ws.Range("B23:B27").Merge()
$ws.Range("D23:D27").Merge()
$ws.Range("F23:F27").Merge()
$ws.Range("B28:B32").Merge()
$ws.Range("E28:E32").Merge()
$ws.Range("C6:C18").Merge()

# ------------------------------------------------------------------
# 4) The multi-line values above trigger an automatic custom row
#    height on their row; AutoFit puts each row back to the sheet's
#    standard (non-custom) height.
# ------------------------------------------------------------------
$ws.Rows(6).AutoFit()
$ws.Rows(23).AutoFit()
$ws.Rows(28).AutoFit()

# ------------------------------------------------------------------
# 5) Remove the old GS-QC-6301 course block (old B6:B12 / F6:F12
#    merge). Unmerge, drop back to the Normal style (so no leftover
#    styled-but-empty cell survives) and clear the content.
# ------------------------------------------------------------------
$ws.Range("B6:B12").UnMerge()
$ws.Range("F6:F12").UnMerge()
$ws.Range("B6:B12").Style = "Normal"
$ws.Range("F6:F12").Style = "Normal"
$ws.Range("B6:B12").ClearContents()
$ws.Range("F6:F12").ClearContents()

Write-Host "edit complete"
